# Add a "cfop" column to the "PI hours" sheet, and add a new "cfop hours" sheet
# summarizing cfop totals, mirroring the structure of the existing
# dept/app summary sheets.

$wb = $excel.ActiveWorkbook

$piSheet = $wb.Worksheets.Item("PI hours")

# --- 1. Add the new "cfop" column (G) to the "PI hours" sheet ---
$piSheet.Range("E1").Copy()
$piSheet.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$piSheet.Range("G1").Value = "cfop"
$piSheet.Range("G2").Value = "['cfop_NH']"
$piSheet.Range("G3").Value = "['cfop_SELIG']"

# --- 2. Add the new "cfop hours" sheet after "unit(accumulative) hours" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cfopSheet = $wb.Worksheets.Add($null, $lastSheet)
$cfopSheet.Name = "cfop hours"

# copy header/index formatting from the "PI hours" sheet
$piSheet.Range("E1").Copy()
$cfopSheet.Range("B1:D1").PasteSpecial(-4122)
$piSheet.Range("A2").Copy()
$cfopSheet.Range("A2").PasteSpecial(-4122)
$cfopSheet.Range("A3").PasteSpecial(-4122)

$cfopSheet.Range("B1").Value = "cfop"
$cfopSheet.Range("C1").Value = "hours"
$cfopSheet.Range("D1").Value = "percentage"

$cfopSheet.Range("A2").Value = 0
$cfopSheet.Range("B2").Value = "cfop_NH"
$cfopSheet.Range("C2").Value = 61
$cfopSheet.Range("D2").Value = 88.40579710144928

$cfopSheet.Range("A3").Value = 1
$cfopSheet.Range("B3").Value = "cfop_SELIG"
$cfopSheet.Range("C3").Value = 8
$cfopSheet.Range("D3").Value = 11.59420289855072

$piSheet.Select()

Write-Host "cfop column and cfop hours sheet added"
